$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.475.39"
$ws.Range("E2").Value = "  -1.14%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.105.83"
$ws.Range("E3").Value = "  -0.51%  "

$ws.Range("E4").Value = "  +0.26%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.12"
$ws.Range("E5").Value = "  +0.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5255"
$ws.Range("E7").Value = "  -1.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4588"
$ws.Range("E8").Value = "  +3.86%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.20"
$ws.Range("E9").Value = "  +12.55%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08993"
$ws.Range("E10").Value = "  -0.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.180"
$ws.Range("E11").Value = "  +0.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.36"
$ws.Range("E12").Value = "  -2.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.100.02"
$ws.Range("E13").Value = "  -0.74%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.785"
$ws.Range("E14").Value = "  +0.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.837"
$ws.Range("E15").Value = "  +0.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.59"
$ws.Range("E16").Value = "  -0.26%  "

$ws.Range("E17").Value = "  +0.30%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001129"
$ws.Range("E18").Value = "  -0.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06627"
$ws.Range("E19").Value = "  -0.80%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.49"
$ws.Range("E20").Value = "  +1.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.318"
$ws.Range("E22").Value = "  -0.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.526.22"
$ws.Range("E23").Value = "  -1.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.35"
$ws.Range("E24").Value = "  +0.53%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.352"
$ws.Range("E25").Value = "  +3.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.357.09"
$ws.Range("E26").Value = "  -0.26%  "

$ws.Range("E27").Value = "  -1.79%  "

$ws.Range("E28").Value = "  -1.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.89"
$ws.Range("E29").Value = "  +0.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.04"
$ws.Range("E30").Value = "  -0.31%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.194"
$ws.Range("E31").Value = "  -0.18%  "

$ws.Range("E32").Value = "  -1.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.660"
$ws.Range("E33").Value = "  +5.82%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.153"
$ws.Range("E34").Value = "  -1.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.931"
$ws.Range("E35").Value = "  -2.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.50"
$ws.Range("E36").Value = "  +9.65%  "

$ws.Range("E37").Value = "  -0.86%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06847"
$ws.Range("E38").Value = "  +0.72%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.544"
$ws.Range("E39").Value = "  -0.18%  "

$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2299"
$ws.Range("E40").Value = "  -0.26%  "

$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.74"
$ws.Range("E41").Value = "  -1.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6901"
$ws.Range("E42").Value = "  +0.76%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.247"
$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.354"
$ws.Range("E44").Value = "  +4.88%  "

$ws.Range("E45").Value = "  +0.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "14.04"
$ws.Range("E46").Value = "  -0.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6386"
$ws.Range("E47").Value = "  -0.77%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.659"
$ws.Range("E48").Value = "  +0.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.245"
$ws.Range("E49").Value = "  -1.85%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.224"
$ws.Range("E50").Value = "  +0.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "83.44"
$ws.Range("E51").Value = "  +0.51%  "
